$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, matching style of existing header row (s="1")
$ws.Range("G1").Value = "MSE_median"
$ws.Range("H1").Value = "MAE_median"
$ws.Range("I1").Value = "Dir_accuracy"

# Copy style from an existing header cell (F1) to the new header cells
$ws.Range("F1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2 new values
$ws.Range("G2").Value = 0.0005099818166665428
$ws.Range("H2").Value = 0.02258274386707427
$ws.Range("I2").Value = 0.4794326241134752

# Row 3 new values
$ws.Range("G3").Value = 0.001083194070471167
$ws.Range("H3").Value = 0.03291191380748265
$ws.Range("I3").Value = ""
